# Weekly update: a new price-survey record (week of 2023-08-28) is inserted
# as a new row 619 in the "Pepino ensalada" sheet. All the rows that used to
# live at 619..669 shift down by one (to 620..670); the previously-last row
# (old 669) ends up as the new row 670.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 619..669 down to 620..670, leaving a blank row 619
# (this also grows the sheet dimension from A1:R669 to A1:R670).
$ws.Rows(619).Insert()

# Populate the newly inserted row 619 with the new survey record.
$ws.Range("A619").Value = 6
$ws.Range("B619").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C619").Value = "Metropolitana"
$ws.Range("D619").Value = 45166
$ws.Range("D619").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E619").Value = 13
$ws.Range("F619").Value = 100112043
$ws.Range("G619").Value = "Pepino ensalada"
$ws.Range("H619").Value = "Sin especificar"
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 570
$ws.Range("K619").Value = 6000
$ws.Range("L619").Value = 7000
$ws.Range("M619").Value = 6439
$ws.Range("N619").Value = "`$/caja 60 unidades"
$ws.Range("O619").Value = "Región de Arica y Parinacota"
$ws.Range("P619").Value = 107
$ws.Range("Q619").Value = 60
$ws.Range("R619").Value = "Hortaliza"
